# ---------------------------------------------------------------------------
# Add a new "2022-Q1" sheet (positioned between "2021-Q4" and "总计"),
# populate it with the quarterly fund-holding breakdown, and append the
# matching summary row to the "总计" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Create the new sheet by duplicating "2021-Q4" (same columns/styles)
$template = $wb.Worksheets.Item("2021-Q4")
$total    = $wb.Worksheets.Item("总计")
$template.Copy($total)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# "总计" shifted one slot to the right when the copy was inserted in front of
# it, and this runtime's Worksheet handles are index-bound - re-resolve the
# reference by name now so later calls hit the right sheet.
$total = $wb.Worksheets.Item("总计")

# ---- 2. Wipe the copied body rows' VALUES only (keep the header row intact,
#         and keep column-A's numbered-index style on rows 2-11 so it doesn't
#         have to be re-applied below).
$usedRows = $q1.UsedRange.Rows.Count
if ($usedRows -gt 1) {
    $q1.Range("A2:H" + $usedRows).ClearContents()
}

# ---- 3. Row data: index, fund code, fund name, fund size, stock position,
#         position ratio, holding value (billion yuan), position rank
$q1Data = @(
    @('0', '010885', '长盛优势企业精选混合A', '15.37', '79.73', '3.23', '0.4965', '8'),
    @('1', '160813', '长盛同盛成长优选灵活配置混合 (LOF)', '12.92', '77.90', '2.93', '0.3786', '8'),
    @('2', '519039', '长盛同德主题混合', '12.19', '79.80', '3.02', '0.3681', '9'),
    @('3', '005396', '中金丰硕混合', '1.83', '71.47', '6.62', '0.1211', '4'),
    @('4', '000534', '长盛高端装备制造灵活配置混合', '3.03', '74.39', '3.45', '0.1045', '4'),
    @('5', '010155', '长盛核心成长混合A', '3.09', '75.19', '2.95', '0.0912', '9'),
    @('6', '006279', '中金瑞祥灵活配置混合A', '2.10', '59.54', '4.03', '0.0846', '9'),
    @('7', '000598', '长盛生态环境主题灵活配置混合', '1.70', '80.08', '3.91', '0.0665', '4'),
    @('8', '001370', '中银新趋势灵活配置混合', '2.63', '34.29', '2.45', '0.0644', '4'),
    @('9', '007305', '国联安新科技混合', '2.13', '81.55', '2.47', '0.0526', '9'),
    @('10', '002156', '长盛盛世灵活配置混合A', '3.20', '28.47', '1.15', '0.0368', '8'),
    @('11', '080002', '长盛创新先锋混合', '0.83', '74.24', '3.65', '0.0303', '5'),
    @('12', '004332', '恒生前海沪港深新兴产业精选混合', '0.52', '80.98', '4.98', '0.0259', '1'),
    @('13', '010886', '长盛优势企业精选混合C', '0.79', '79.73', '3.23', '0.0255', '8'),
    @('14', '005728', '华宝绿色主题混合', '0.49', '82.30', '4.10', '0.0201', '8'),
    @('15', '005903', '泰达宏利绩优增长灵活配置混合', '0.40', '92.75', '4.34', '0.0174', '5'),
    @('16', '010156', '长盛核心成长混合C', '0.37', '75.19', '2.95', '0.0109', '9'),
    @('17', '002157', '长盛盛世灵活配置混合C', '0.20', '28.47', '1.15', '0.0023', '8'),
    @('18', '006280', '中金瑞祥灵活配置混合C', '0.00', '59.54', '4.03', '0', '9')
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percentages /
# values with significant trailing zeros) that must stay text, exactly like
# the sibling quarterly sheets. Pre-format the destination range as Text so
# values such as "010885" or "79.73" are not silently coerced to numbers.
$lastRow = 1 + $q1Data.Count
$q1.Range("B2:B" + $lastRow).NumberFormat = "@"
$q1.Range("D2:G" + $lastRow).NumberFormat = "@"

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $r = 2 + $i
    $row = $q1Data[$i]

    $q1.Cells.Item($r, 1).Value = [double]$row[0]          # A - row index (number)
    $q1.Cells.Item($r, 2).Value = $row[1]                  # B - fund code (text)
    $q1.Cells.Item($r, 3).Value = $row[2]                  # C - fund name (text)
    $q1.Cells.Item($r, 4).Value = $row[3]                  # D - fund size (text)
    $q1.Cells.Item($r, 5).Value = $row[4]                  # E - stock position (text)
    $q1.Cells.Item($r, 6).Value = $row[5]                  # F - position ratio (text)

    # G - holding value: numeric zero stays a real number, everything else
    # keeps the text formatting applied above (matches "2021-Q4" row 11).
    if ($row[6] -eq "0") {
        $q1.Range("G" + $r).NumberFormat = "General"
        $q1.Cells.Item($r, 7).Value = 0
    } else {
        $q1.Cells.Item($r, 7).Value = $row[6]
    }

    $q1.Cells.Item($r, 8).Value = [double]$row[7]          # H - position rank (number)
}

# ---- 4. Re-apply the first-column "index" style to the rows beyond the
#         11 originally copied from "2021-Q4" (rows 12-20 here).
$template.Range("A11").Copy()
if ($lastRow -gt 11) {
    $q1.Range("A12:A" + $lastRow).PasteSpecial(-4122)
}
$q1.Range("A1").Select()

# ---- 5. Insert the 2022-Q1 summary row at the top of the "总计" sheet body
$total.Rows(2).Insert()
$total.Rows(2).ClearFormats()
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 19
$total.Cells.Item(2, 4).Value = 2

# Column A keeps the same "index" style used by every other row in the sheet.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A1").Select()
